$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# ---------------------------------------------------------------------------
# 1) Re-shuffle the formatting of the last three table rows BEFORE inserting,
#    so that after the insert the newly added last row inherits the old
#    "final row" bottom-border style, the old final row becomes a normal
#    in-table row, and the row above it takes the "second to last" style.
# ---------------------------------------------------------------------------

# Capture the current last data row's (row 138) formatting - this becomes
# the formatting of the brand new row (139) once the table grows.
$ws.Range("A138:K138").Copy() | Out-Null

# Grow the table by one row (extends Table1 ref to A8:K139, sheet dimension
# to A2:K139).
$tbl.ListRows.Add() | Out-Null

# Paste the captured "final row" formatting onto the new last row (139).
$ws.Range("A139:K139").PasteSpecial(-4122) | Out-Null
$ws.Range("G139").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# Row 137 (previously the "final row" style) now becomes the "second to
# last" row - copy what used to be row 137's own formatting down into 138.
$ws.Range("A137:K137").Copy() | Out-Null
$ws.Range("A138:K138").PasteSpecial(-4122) | Out-Null
$ws.Range("G138").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# Row 137 itself becomes an ordinary in-table row, matching row 136's style.
$ws.Range("A136:K136").Copy() | Out-Null
$ws.Range("A137:K137").PasteSpecial(-4122) | Out-Null
$ws.Range("G137").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) New leave-card entries for FY2023 (Aug-Dec) and the 2024 header/rows.
#    New shared strings must be introduced in the same order they appear in
#    the authored workbook: "2024" (row 34) before "SL(5-0-0)" (row 31)
#    before "10/16,23-26/2023" (row 31).
# ---------------------------------------------------------------------------

# Row 34 - new "2024" year header (matches the style used by the 2022/2023
# headers in rows 10 and 18). Do this first so the "2024" shared string is
# allocated before the other two new strings below.
$ws.Range("A18:K18").Copy() | Out-Null
$ws.Range("A34:K34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("A34").Value = "'2024"

# Row 29 - August 2023: record 1.25 VL/SL earned this period.
$ws.Range("C29").Value = 1.25

# Row 30 - September 2023.
$ws.Range("A30").Value = 45170
$ws.Range("C30").Value = 1.25

# Row 31 - October 2023, sick leave usage entry.
$ws.Range("A31").Value = 45200
$ws.Range("B31").Value = "SL(5-0-0)"
$ws.Range("C31").Value = 1.25
$ws.Range("H31").Value = 5
$ws.Range("K31").Value = "10/16,23-26/2023"

# Row 32 - November 2023, vacation leave usage entry.
$ws.Range("A32").Value = 45231
$ws.Range("B32").Value = "VL(1-0-0)"
$ws.Range("D32").Value = 1
$ws.Range("K12").Copy() | Out-Null
$ws.Range("K32").PasteSpecial(-4122) | Out-Null
$ws.Range("K32").Value = 45261
$excel.CutCopyMode = $false

# Row 33 - December 2023.
$ws.Range("A33").Value = 45261

# Rows 35-42 - first-of-month period markers for 2024.
$ws.Range("A35").Value = 45292
$ws.Range("A36").Value = 45323
$ws.Range("A37").Value = 45352
$ws.Range("A38").Value = 45383
$ws.Range("A39").Value = 45413
$ws.Range("A40").Value = 45444
$ws.Range("A41").Value = 45474
$ws.Range("A42").Value = 45505

# ---------------------------------------------------------------------------
# 3) Leave current selection positioned like the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("K32").Select() | Out-Null

$wb.Save()
